# Bugfix to satisfaction criteria settings
#
# The "Satisfaction criterion mode" column (F) on the "Internal drivers"
# sheet was incorrectly set to "Minimise" for every driver row; it should
# read "Maximise" (a new shared string the workbook didn't have before).
# Only column F is affected - the other "mode" columns (G:J) legitimately
# stay as "Minimise".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Internal drivers")

# Data rows run from 3 to 34 (rows 1-2 are header rows).
$ws.Range("F3:F34").Value = "Maximise"

# Leave the cursor where the author ended up after making the fix.
$ws.Range("I13").Select()
